# Generate Report for handoff
# f26dc623-...md moves to row2 (status -> "In Translation"); b67cb069-...md moves
# to row3 (status stays "Ready for handoff" but gets a fresh handoff datetime).

$wb = $excel.ActiveWorkbook

# ---------- Overview sheet ----------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "f26dc623-f063-4b95-9d3c-845b0b07d230.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "In Translation"
$ws.Range("A3").Value = "b67cb069-9920-4cfe-adfc-5112adfff8ff.md"

# refresh hyperlinks so display text tracks the swapped file names
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/02419727e64cdff9e9318b37ee9ab4133019ba1b/e2e/f26dc623-f063-4b95-9d3c-845b0b07d230.md", "", "", "f26dc623-f063-4b95-9d3c-845b0b07d230.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/02419727e64cdff9e9318b37ee9ab4133019ba1b/e2e/b67cb069-9920-4cfe-adfc-5112adfff8ff.md", "", "", "b67cb069-9920-4cfe-adfc-5112adfff8ff.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/02419727e64cdff9e9318b37ee9ab4133019ba1b/.localization-config", "", "", ".localization-config") | Out-Null

# ---------- zh-cn sheet ----------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = "f26dc623-f063-4b95-9d3c-845b0b07d230.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "f26dc623-f063-4b95-9d3c-845b0b07d230.07a9b1d56c24218aed9edbcf5d91573c4bb8b627.zh-cn.xlf"
$ws.Range("A3").Value = "b67cb069-9920-4cfe-adfc-5112adfff8ff.md"
$ws.Range("C3").Value = "b67cb069-9920-4cfe-adfc-5112adfff8ff.3f88f1fb3b7fc503e96438ffd2cc252740a74de8.zh-cn.xlf"
$ws.Range("D3").Value = "2016-01-26 08:44:51"

$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/02419727e64cdff9e9318b37ee9ab4133019ba1b/e2e/f26dc623-f063-4b95-9d3c-845b0b07d230.md", "", "", "f26dc623-f063-4b95-9d3c-845b0b07d230.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/457df19dc110990d4583ca627fa5ef86dccab499/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f26dc623-f063-4b95-9d3c-845b0b07d230.07a9b1d56c24218aed9edbcf5d91573c4bb8b627.zh-cn.xlf", "", "", "f26dc623-f063-4b95-9d3c-845b0b07d230.07a9b1d56c24218aed9edbcf5d91573c4bb8b627.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/02419727e64cdff9e9318b37ee9ab4133019ba1b/e2e/b67cb069-9920-4cfe-adfc-5112adfff8ff.md", "", "", "b67cb069-9920-4cfe-adfc-5112adfff8ff.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/457df19dc110990d4583ca627fa5ef86dccab499/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/b67cb069-9920-4cfe-adfc-5112adfff8ff.3f88f1fb3b7fc503e96438ffd2cc252740a74de8.zh-cn.xlf", "", "", "b67cb069-9920-4cfe-adfc-5112adfff8ff.3f88f1fb3b7fc503e96438ffd2cc252740a74de8.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/02419727e64cdff9e9318b37ee9ab4133019ba1b/.localization-config", "", "", ".localization-config") | Out-Null

# ---------- de-de sheet ----------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = "f26dc623-f063-4b95-9d3c-845b0b07d230.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "f26dc623-f063-4b95-9d3c-845b0b07d230.07a9b1d56c24218aed9edbcf5d91573c4bb8b627.de-de.xlf"
$ws.Range("A3").Value = "b67cb069-9920-4cfe-adfc-5112adfff8ff.md"
$ws.Range("C3").Value = "b67cb069-9920-4cfe-adfc-5112adfff8ff.3f88f1fb3b7fc503e96438ffd2cc252740a74de8.de-de.xlf"
$ws.Range("D3").Value = "2016-01-26 08:45:03"

$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/02419727e64cdff9e9318b37ee9ab4133019ba1b/e2e/f26dc623-f063-4b95-9d3c-845b0b07d230.md", "", "", "f26dc623-f063-4b95-9d3c-845b0b07d230.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3e7f28cf4a8464b9ffa1b0a2342be5f70c4fd9bd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f26dc623-f063-4b95-9d3c-845b0b07d230.07a9b1d56c24218aed9edbcf5d91573c4bb8b627.de-de.xlf", "", "", "f26dc623-f063-4b95-9d3c-845b0b07d230.07a9b1d56c24218aed9edbcf5d91573c4bb8b627.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/02419727e64cdff9e9318b37ee9ab4133019ba1b/e2e/b67cb069-9920-4cfe-adfc-5112adfff8ff.md", "", "", "b67cb069-9920-4cfe-adfc-5112adfff8ff.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3e7f28cf4a8464b9ffa1b0a2342be5f70c4fd9bd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/b67cb069-9920-4cfe-adfc-5112adfff8ff.3f88f1fb3b7fc503e96438ffd2cc252740a74de8.de-de.xlf", "", "", "b67cb069-9920-4cfe-adfc-5112adfff8ff.3f88f1fb3b7fc503e96438ffd2cc252740a74de8.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/02419727e64cdff9e9318b37ee9ab4133019ba1b/.localization-config", "", "", ".localization-config") | Out-Null
